$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.0599
$ws.Range("A3").Value = -21.37310000000003
$ws.Range("E5").Value = 13.1298
$ws.Range("A14").Value = -20.57909999999998
$ws.Range("A21").Value = -21.28120000000001
$ws.Range("A23").Value = -21.40160000000002
$ws.Range("A25").Value = -22.43880000000003
